$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying the format of the existing
# header cells (e.g. G1 "sum") so it matches the other headers
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the Save column values for the data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
